$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 530.7143
$ws.Range("I9").Value = 463
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 463
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -294
$ws.Range("N9").Value = -1038
$ws.Range("H18").Value = 720.5
$ws.Range("I18").Value = 720.5
$ws.Range("K18").Value = 720.5
$ws.Range("M18").Value = -436.5
$ws.Range("H28").Value = 1750
$ws.Range("I28").Value = 1298.8
$ws.Range("J28").Value = 4006
$ws.Range("K28").Value = 1298.8
$ws.Range("L28").Value = 4006
$ws.Range("M28").Value = -813.8
$ws.Range("N28").Value = -4976
$ws.Range("H80").Value = 3403.9333
$ws.Range("I80").Value = 663.6667
$ws.Range("J80").Value = 4089
$ws.Range("K80").Value = 1991.0001
$ws.Range("L80").Value = 12267
$ws.Range("M80").Value = -993.0001
$ws.Range("N80").Value = -14263
$ws.Range("H83").Value = 3403.9333
$ws.Range("I83").Value = 663.6667
$ws.Range("J83").Value = 4089
$ws.Range("K83").Value = 5973.0003
$ws.Range("L83").Value = 36801
$ws.Range("M83").Value = -981.0002999999997
$ws.Range("N83").Value = -46785
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H137").Value = 14651.305
$ws.Range("I137").Value = 3180.926
$ws.Range("K137").Value = 9542.778
$ws.Range("M137").Value = -6992.778
$ws.Range("H138").Value = 1472.5294
$ws.Range("I138").Value = 1283.2
$ws.Range("J138").Value = 2892.5
$ws.Range("K138").Value = 3849.6
$ws.Range("L138").Value = 8677.5
$ws.Range("M138").Value = 1290.4
$ws.Range("N138").Value = -18957.5
$ws.Range("N106").ClearContents()

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1995.0714
$ws.Range("I2").Value = 923
$ws.Range("J2").Value = 3924.8
$ws.Range("K2").Value = 923
$ws.Range("L2").Value = 3924.8
$ws.Range("M2").Value = -810
$ws.Range("N2").Value = -4150.8
$ws.Range("H22").Value = 3250
$ws.Range("I22").Value = 3250
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3250
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2951
$ws.Range("H31").Value = 12466.333
$ws.Range("I31").Value = 12466.333
$ws.Range("K31").Value = 12466.333
$ws.Range("M31").Value = -12172.333
$ws.Range("H32").Value = 1405.1333
$ws.Range("I32").Value = 1405.1333
$ws.Range("K32").Value = 1405.1333
$ws.Range("M32").Value = -1118.1333
$ws.Range("H116").Value = 1995.0714
$ws.Range("I116").Value = 923
$ws.Range("J116").Value = 3924.8
$ws.Range("K116").Value = 923
$ws.Range("L116").Value = 3924.8
$ws.Range("M116").Value = 1371
$ws.Range("N116").Value = -8512.799999999999
$ws.Range("H122").Value = 1814.3334
$ws.Range("I122").Value = 1814.3334
$ws.Range("K122").Value = 5443.0002
$ws.Range("M122").Value = -2993.0002
$ws.Range("N22").ClearContents()

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1995.0714
$ws.Range("I3").Value = 923
$ws.Range("J3").Value = 3924.8
$ws.Range("K3").Value = 923
$ws.Range("L3").Value = 3924.8
$ws.Range("M3").Value = -809
$ws.Range("N3").Value = -4152.8
$ws.Range("H18").Value = 52000
$ws.Range("J18").Value = 52000
$ws.Range("L18").Value = 52000
$ws.Range("N18").Value = -53058
$ws.Range("H20").Value = 3200.6
$ws.Range("I20").Value = 2001.3334
$ws.Range("K20").Value = 2001.3334
$ws.Range("M20").Value = -1754.3334
$ws.Range("H99").Value = 6079.7144

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18544.834
$ws.Range("I31").Value = 32459.812
$ws.Range("K31").Value = 32459.812
$ws.Range("M31").Value = -32164.812
$ws.Range("H34").Value = 18544.834
$ws.Range("I34").Value = 32459.812
$ws.Range("K34").Value = 32459.812
$ws.Range("M34").Value = -32257.812
$ws.Range("H62").Value = 34779
$ws.Range("I62").Value = 29998
$ws.Range("K62").Value = 29998
$ws.Range("M62").Value = -29374
$ws.Range("H65").Value = 34779
$ws.Range("I65").Value = 29998
$ws.Range("K65").Value = 149990
$ws.Range("M65").Value = -146870
$ws.Range("H107").Value = 4244
$ws.Range("J107").Value = 4007.2
$ws.Range("L107").Value = 4007.2
$ws.Range("N107").Value = -7847.2
$ws.Range("H132").Value = 5379.3687
$ws.Range("I132").Value = 4894.5884
$ws.Range("K132").Value = 14683.7652
$ws.Range("M132").Value = -12153.7652
$ws.Range("H134").Value = 7145430
$ws.Range("I134").Value = 2784.818
$ws.Range("K134").Value = 8354.454000000002
$ws.Range("M134").Value = -5819.454000000002

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 830.625
$ws.Range("I8").Value = 830.625
$ws.Range("K8").Value = 2491.875
$ws.Range("M8").Value = -2352.875
$ws.Range("H12").Value = 1250.7273
$ws.Range("J12").Value = 1531
$ws.Range("L12").Value = 4593
$ws.Range("N12").Value = -4939
$ws.Range("H37").Value = 59999
$ws.Range("J37").Value = 59999
$ws.Range("L37").Value = 179997
$ws.Range("N37").Value = -180221
$ws.Range("H47").Value = 1761.5333
$ws.Range("I47").Value = 119.28571
$ws.Range("J47").Value = 3198.5
$ws.Range("K47").Value = 357.85713
$ws.Range("L47").Value = 9595.5
$ws.Range("M47").Value = 73.14287000000002
$ws.Range("N47").Value = -10457.5
$ws.Range("H99").Value = 3517.25
$ws.Range("I99").Value = 3517.25
$ws.Range("K99").Value = 10551.75
$ws.Range("M99").Value = -8305.75
$ws.Range("H131").Value = 1213271.5
$ws.Range("I131").Value = 881.875
$ws.Range("K131").Value = 2645.625
$ws.Range("M131").Value = 2394.375
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 964.8
$ws.Range("J132").Value = 1719
$ws.Range("K132").Value = 8683.199999999999
$ws.Range("L132").Value = 15471
$ws.Range("M132").Value = -6153.199999999999
$ws.Range("N132").Value = -20531
$ws.Range("H140").Value = 1738
$ws.Range("I140").Value = 1383.4546
$ws.Range("J140").Value = 2518
$ws.Range("K140").Value = 4150.3638
$ws.Range("L140").Value = 7554
$ws.Range("M140").Value = 1029.6362
$ws.Range("N140").Value = -17914

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1998.3334
$ws.Range("J102").Value = 1998.25
$ws.Range("L102").Value = 1998.25
$ws.Range("N102").Value = -5242.25
$ws.Range("H126").Value = 4279.722
$ws.Range("I126").Value = 5899
$ws.Range("J126").Value = 2984.3
$ws.Range("K126").Value = 17697
$ws.Range("L126").Value = 8952.900000000001
$ws.Range("M126").Value = -15227
$ws.Range("N126").Value = -13892.9
$ws.Range("H132").Value = 2440.3333
$ws.Range("I132").Value = 2339
$ws.Range("J132").Value = 2947
$ws.Range("K132").Value = 7017
$ws.Range("L132").Value = 8841
$ws.Range("M132").Value = -4487
$ws.Range("N132").Value = -13901

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2269.4285
$ws.Range("J22").Value = 2235.75
$ws.Range("L22").Value = 2235.75
$ws.Range("N22").Value = -2825.75
$ws.Range("H27").Value = 2269.4285
$ws.Range("J27").Value = 2235.75
$ws.Range("L27").Value = 2235.75
$ws.Range("N27").Value = -2449.75
$ws.Range("H46").Value = 2442.4348
$ws.Range("I46").Value = 1639.4
$ws.Range("J46").Value = 2665.5
$ws.Range("K46").Value = 1639.4
$ws.Range("L46").Value = 2665.5
$ws.Range("M46").Value = -1451.4
$ws.Range("N46").Value = -3041.5
$ws.Range("H55").Value = 1527.3334
$ws.Range("I55").Value = 293.5
$ws.Range("K55").Value = 293.5
$ws.Range("M55").Value = -120.5
$ws.Range("H82").Value = 3200.1428
$ws.Range("J82").Value = 8751.5
$ws.Range("L82").Value = 8751.5
$ws.Range("N82").Value = -9473.5
$ws.Range("H85").Value = 3200.1428
$ws.Range("J85").Value = 8751.5
$ws.Range("L85").Value = 8751.5
$ws.Range("N85").Value = -11247.5
$ws.Range("H122").Value = 3056.8845
$ws.Range("I122").Value = 2680.6924
$ws.Range("J122").Value = 3433.077
$ws.Range("K122").Value = 8042.0772
$ws.Range("L122").Value = 10299.231
$ws.Range("M122").Value = -5592.0772
$ws.Range("N122").Value = -15199.231
$ws.Range("H136").Value = 2149.6667
$ws.Range("I136").Value = 1654.1818
$ws.Range("K136").Value = 4962.5454
$ws.Range("M136").Value = -2412.5454

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9388.388999999999
$ws.Range("I81").Value = 10768.538
$ws.Range("K81").Value = 21537.076
$ws.Range("M81").Value = -20476.076
$ws.Range("H84").Value = 9388.388999999999
$ws.Range("I84").Value = 10768.538
$ws.Range("K84").Value = 107685.38
$ws.Range("M84").Value = -102381.38
$ws.Range("H107").Value = 1196.4348
$ws.Range("I107").Value = 734.8889
$ws.Range("K107").Value = 2204.6667
$ws.Range("M107").Value = -284.6667000000002
$ws.Range("H118").Value = 30106.4
$ws.Range("I118").Value = 29356
$ws.Range("J118").Value = 30294
$ws.Range("K118").Value = 29356
$ws.Range("L118").Value = 30294
$ws.Range("M118").Value = -27699
$ws.Range("N118").Value = -33608
$ws.Range("H126").Value = 4799.3335
$ws.Range("I126").Value = 4249.5
$ws.Range("K126").Value = 12748.5
$ws.Range("M126").Value = -10278.5
$ws.Range("H132").Value = 2256.5312
$ws.Range("I132").Value = 2280.9355
$ws.Range("K132").Value = 6842.806500000001
$ws.Range("M132").Value = -4312.806500000001
$ws.Range("H136").Value = 1841.2
$ws.Range("I136").Value = 1758.5
$ws.Range("K136").Value = 5275.5
$ws.Range("M136").Value = -2725.5
